$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("state_id")

# Rename the column header in B1 from "stcode11" to "code11"
$ws.Range("B1").Value = "code11"

# Move the active selection to B1 (matches post-edit selection state)
$ws.Range("B1").Select()
